$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing row (76) down onto the new row (77)
# so the new row inherits the same cell styles as the rest of the table.
$xlPasteFormats = -4122
$ws.Range("A76:H76").Copy()
$ws.Range("A77:H77").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# Fill in the new journal entry (row 77)
$ws.Range("A77").Value2 = 44638
$ws.Range("B77").Value2 = 0.59722222222222221
$ws.Range("C77").Value2 = 0.62430555555555556
$ws.Range("D77").Formula = "=Tableau4[[#This Row],[Heure fin]]-Tableau4[[#This Row],[Heure début]]"
$ws.Range("E77").Value2 = "CPNV"
$ws.Range("F77").Value2 = "Créer les 5 grilles"
$ws.Range("G77").Value2 = "M. Viret m'a expliqué et aidé sur Teams comment lire un fichier. J'ai donc plus de problèmes"
$ws.Range("H77").Value2 = "M. Viret"

# Grow the "Tableau4" table range to include the new row
$lo = $ws.ListObjects.Item("Tableau4")
$lo.Resize($ws.Range("A1:H77"))

# Move the active selection to the new last cell, matching the author's cursor position
$ws.Range("H77").Select()
